$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 242 - this shifts the existing rows
# 242-320 down to 243-321 (matches the diff's row-shift pattern) and
# extends the used range to A1:T321.
$ws.Rows.Item(242).Insert()

# Populate the newly inserted row 242 with the new weekly record.
# Columns that are identical across all "Mango" rows (A,B,C,E,F,G,H,I,J,K,Q,T)
# are copied from the template; L (Calidad) and R (Origen) are unchanged
# from the row that used to occupy position 242; D,M,N,O,P,S are the new
# values introduced by this edit.
$ws.Cells.Item(242, 1).Value = 9
$ws.Cells.Item(242, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(242, 3).Value = "Metropolitana"
$ws.Cells.Item(242, 4).Value = 44559
$ws.Cells.Item(242, 5).Value = 13
$ws.Cells.Item(242, 6).Value = "Fruta"
$ws.Cells.Item(242, 7).Value = 100108
$ws.Cells.Item(242, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(242, 9).Value = 100108002
$ws.Cells.Item(242, 10).Value = "Mango"
$ws.Cells.Item(242, 11).Value = "Sin especificar"
$ws.Cells.Item(242, 12).Value = "Primera"
$ws.Cells.Item(242, 13).Value = 720
$ws.Cells.Item(242, 14).Value = 5000
$ws.Cells.Item(242, 15).Value = 5500
$ws.Cells.Item(242, 16).Value = 5243
$ws.Cells.Item(242, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(242, 18).Value = "Perú"
$ws.Cells.Item(242, 19).Value = 1311
$ws.Cells.Item(242, 20).Value = 4
